$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 6444.778
$ws.Range("I76").Value = 8500.75
$ws.Range("J76").Value = 4800
$ws.Range("K76").Value = 8500.75
$ws.Range("L76").Value = 4800
$ws.Range("M76").Value = -8185.75
$ws.Range("N76").Value = -5430

$ws.Range("H79").Value = 6444.778
$ws.Range("I79").Value = 8500.75
$ws.Range("J79").Value = 4800
$ws.Range("K79").Value = 8500.75
$ws.Range("L79").Value = 4800
$ws.Range("M79").Value = -7408.75
$ws.Range("N79").Value = -6984

$ws.Range("H86").Value = 1986.1428
$ws.Range("I86").Value = 2150.75
$ws.Range("J86").Value = 1766.6666
$ws.Range("K86").Value = 2150.75
$ws.Range("L86").Value = 1766.6666
$ws.Range("M86").Value = -1027.75
$ws.Range("N86").Value = -4012.6666

$ws.Range("H89").Value = 1986.1428
$ws.Range("I89").Value = 2150.75
$ws.Range("J89").Value = 1766.6666
$ws.Range("K89").Value = 10753.75
$ws.Range("L89").Value = 8833.333000000001
$ws.Range("M89").Value = -5137.75
$ws.Range("N89").Value = -20065.333

$ws.Range("H129").Value = 1060.4082
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 1060.4082
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 3181.2246
$ws.Range("N129").Value = -13181.2246

$ws.Range("H137").Value = 2372.625
$ws.Range("I137").Value = 1729.6471
$ws.Range("J137").Value = 2847.8696
$ws.Range("K137").Value = 5188.9413
$ws.Range("L137").Value = 8543.6088
$ws.Range("M137").Value = -2638.9413
$ws.Range("N137").Value = -13643.6088

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("M56").ClearContents()

$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()

$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("M60").ClearContents()

$ws.Range("H61").Value = 2807.4092
$ws.Range("I61").Value = 2750.8667
$ws.Range("J61").Value = 2928.5715
$ws.Range("K61").Value = 2750.8667
$ws.Range("L61").Value = 2928.5715
$ws.Range("M61").Value = -2538.8667
$ws.Range("N61").Value = -3352.5715

$ws.Range("H63").Value = 3093.7273
$ws.Range("I63").Value = 2802.5
$ws.Range("J63").Value = 6006
$ws.Range("K63").Value = 2802.5
$ws.Range("L63").Value = 6006
$ws.Range("M63").Value = -2116.5
$ws.Range("N63").Value = -7378

$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H66").Value = 3093.7273
$ws.Range("I66").Value = 2802.5
$ws.Range("J66").Value = 6006
$ws.Range("K66").Value = 14012.5
$ws.Range("L66").Value = 30030
$ws.Range("M66").Value = -10580.5
$ws.Range("N66").Value = -36894

$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H74").Value = 1922.9584
$ws.Range("I74").Value = 1876.125
$ws.Range("J74").Value = 2016.625
$ws.Range("K74").Value = 1876.125
$ws.Range("L74").Value = 2016.625
$ws.Range("M74").Value = -1002.125
$ws.Range("N74").Value = -3764.625

$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").ClearContents()

$ws.Range("H77").Value = 1922.9584
$ws.Range("I77").Value = 1876.125
$ws.Range("J77").Value = 2016.625
$ws.Range("K77").Value = 9380.625
$ws.Range("L77").Value = 10083.125
$ws.Range("M77").Value = -5012.625
$ws.Range("N77").Value = -18819.125

$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("M78").ClearContents()

$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("N84").ClearContents()

$ws.Range("H121").Value = 31748.334
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 31748.334
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 31748.334
$ws.Range("N121").Value = -35242.334

$ws.Range("H136").Value = 2807.4092
$ws.Range("I136").Value = 2750.8667
$ws.Range("J136").Value = 2928.5715
$ws.Range("K136").Value = 8252.6001
$ws.Range("L136").Value = 8785.7145
$ws.Range("M136").Value = -5702.6001
$ws.Range("N136").Value = -13885.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 25130.912
$ws.Range("I107").Value = 36593.266
$ws.Range("J107").Value = 3639
$ws.Range("K107").Value = 36593.266
$ws.Range("L107").Value = 3639
$ws.Range("M107").Value = -34673.266
$ws.Range("N107").Value = -7479

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 7524.5
$ws.Range("I4").Value = 99.8
$ws.Range("J4").Value = 9999.4
$ws.Range("K4").Value = 99.8
$ws.Range("L4").Value = 9999.4
$ws.Range("M4").Value = 12.2
$ws.Range("N4").Value = -10223.4

$ws.Range("H31").Value = 2723.756
$ws.Range("I31").Value = 2897.4
$ws.Range("J31").Value = 2558.3809
$ws.Range("K31").Value = 2897.4
$ws.Range("L31").Value = 2558.3809
$ws.Range("M31").Value = -2602.4
$ws.Range("N31").Value = -3148.3809

$ws.Range("H34").Value = 2723.756
$ws.Range("I34").Value = 2897.4
$ws.Range("J34").Value = 2558.3809
$ws.Range("K34").Value = 2897.4
$ws.Range("L34").Value = 2558.3809
$ws.Range("M34").Value = -2695.4
$ws.Range("N34").Value = -2962.3809

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 149200
$ws.Range("I4").Value = 200680
$ws.Range("J4").Value = 20500
$ws.Range("K4").Value = 602040
$ws.Range("L4").Value = 61500
$ws.Range("M4").Value = -601928
$ws.Range("N4").Value = -61724

$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()

$ws.Range("H44").Value = 712.2
$ws.Range("I44").Value = 589.1429000000001
$ws.Range("J44").Value = 999.3333
$ws.Range("K44").Value = 1767.4287
$ws.Range("L44").Value = 2997.9999
$ws.Range("M44").Value = -1369.4287
$ws.Range("N44").Value = -3793.9999

$ws.Range("H68").Value = 185912.12
$ws.Range("I68").Value = 200619.5
$ws.Range("J68").Value = 2070
$ws.Range("K68").Value = 601858.5
$ws.Range("L68").Value = 6210
$ws.Range("M68").Value = -601047.5
$ws.Range("N68").Value = -7832

$ws.Range("H71").Value = 185912.12
$ws.Range("I71").Value = 200619.5
$ws.Range("J71").Value = 2070
$ws.Range("K71").Value = 1805575.5
$ws.Range("L71").Value = 18630
$ws.Range("M71").Value = -1801519.5
$ws.Range("N71").Value = -26742

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 27498.5
$ws.Range("I5").Value = 25002
$ws.Range("J5").Value = 29995
$ws.Range("K5").Value = 25002
$ws.Range("L5").Value = 29995
$ws.Range("M5").Value = -24890

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H116").Value = 85294.5
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 85294.5
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 85294.5
$ws.Range("N116").Value = -94472.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5111.1113
$ws.Range("I62").Value = 6000
$ws.Range("J62").Value = 4666.6665
$ws.Range("K62").Value = 6000
$ws.Range("L62").Value = 4666.6665
$ws.Range("M62").Value = -5376
$ws.Range("N62").Value = -5914.6665

$ws.Range("H65").Value = 5111.1113
$ws.Range("I65").Value = 6000
$ws.Range("J65").Value = 4666.6665
$ws.Range("K65").Value = 30000
$ws.Range("L65").Value = 23333.3325
$ws.Range("M65").Value = -26880
$ws.Range("N65").Value = -29573.3325

$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

$ws.Range("H122").Value = 178573170
$ws.Range("I122").Value = 250001420
$ws.Range("J122").Value = 2497.5
$ws.Range("K122").Value = 750004260
$ws.Range("L122").Value = 7492.5
$ws.Range("M122").Value = -750001810
$ws.Range("N122").Value = -12392.5

$ws.Range("H132").Value = 2376.257
$ws.Range("I132").Value = 2384.2632
$ws.Range("J132").Value = 2366.75
$ws.Range("K132").Value = 7152.7896
$ws.Range("L132").Value = 7100.25
$ws.Range("M132").Value = -4622.7896
$ws.Range("N132").Value = -12160.25
